$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $s = $cell.Style
    $cell.Value = "'" + $val
    $cell.Style = $s
}

Set-TextValue $ws.Cells.Item(2,4) "304.78"
Set-TextValue $ws.Cells.Item(2,5) "0.89%"
Set-TextValue $ws.Cells.Item(2,7) "7"
Set-TextValue $ws.Cells.Item(3,4) "36.03"
Set-TextValue $ws.Cells.Item(3,5) "-3.95%"
Set-TextValue $ws.Cells.Item(3,7) "7"
Set-TextValue $ws.Cells.Item(4,4) "5.123"
Set-TextValue $ws.Cells.Item(4,5) "2.26%"
Set-TextValue $ws.Cells.Item(4,7) "7"
Set-TextValue $ws.Cells.Item(5,4) "0.07862"
Set-TextValue $ws.Cells.Item(5,5) "0.11%"
Set-TextValue $ws.Cells.Item(5,7) "7"
Set-TextValue $ws.Cells.Item(6,4) "2.170"
Set-TextValue $ws.Cells.Item(6,5) "-3.15%"
Set-TextValue $ws.Cells.Item(6,7) "7"
Set-TextValue $ws.Cells.Item(7,4) "7.952"
Set-TextValue $ws.Cells.Item(7,5) "-1.03%"
Set-TextValue $ws.Cells.Item(7,7) "7"
Set-TextValue $ws.Cells.Item(8,2) "GateToken"
Set-TextValue $ws.Cells.Item(8,3) "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Cells.Item(8,4) "4.111"
Set-TextValue $ws.Cells.Item(8,5) "2.26%"
Set-TextValue $ws.Cells.Item(8,7) "7"
Set-TextValue $ws.Cells.Item(9,2) "MXToken"
Set-TextValue $ws.Cells.Item(9,3) "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Cells.Item(9,4) "0.9181"
Set-TextValue $ws.Cells.Item(9,5) "1.04%"
Set-TextValue $ws.Cells.Item(9,7) "7"
Set-TextValue $ws.Cells.Item(10,2) "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Cells.Item(10,3) "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Cells.Item(10,4) "0.09659"
Set-TextValue $ws.Cells.Item(10,5) "3.99%"
Set-TextValue $ws.Cells.Item(10,7) "7"
Set-TextValue $ws.Cells.Item(11,2) "WazirX"
Set-TextValue $ws.Cells.Item(11,3) "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Cells.Item(11,4) "0.1871"
Set-TextValue $ws.Cells.Item(11,5) "-0.58%"
Set-TextValue $ws.Cells.Item(11,7) "7"
Set-TextValue $ws.Cells.Item(12,2) "MandalaExchangeToken"
Set-TextValue $ws.Cells.Item(12,3) "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Cells.Item(12,4) "0.08640"
Set-TextValue $ws.Cells.Item(12,5) "1.79%"
Set-TextValue $ws.Cells.Item(12,7) "7"
Set-TextValue $ws.Cells.Item(13,2) "BitrueCoin"
Set-TextValue $ws.Cells.Item(13,3) "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Cells.Item(13,4) "0.03483"
Set-TextValue $ws.Cells.Item(13,5) "-1.16%"
Set-TextValue $ws.Cells.Item(13,7) "7"
Set-TextValue $ws.Cells.Item(14,2) "BitMartToken"
Set-TextValue $ws.Cells.Item(14,3) "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Cells.Item(14,4) "0.09942"
Set-TextValue $ws.Cells.Item(14,5) "-0.25%"
Set-TextValue $ws.Cells.Item(14,7) "7"
Set-TextValue $ws.Cells.Item(15,2) "BitForexToken"
Set-TextValue $ws.Cells.Item(15,3) "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Cells.Item(15,4) "0.001449"
Set-TextValue $ws.Cells.Item(15,5) "-2.34%"
Set-TextValue $ws.Cells.Item(15,7) "7"
Set-TextValue $ws.Cells.Item(16,2) "TigerCash"
Set-TextValue $ws.Cells.Item(16,3) "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Cells.Item(16,4) "0.005706"
Set-TextValue $ws.Cells.Item(16,5) "0.09%"
Set-TextValue $ws.Cells.Item(16,7) "7"
Set-TextValue $ws.Cells.Item(17,2) "LEO"
Set-TextValue $ws.Cells.Item(17,3) "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Cells.Item(17,4) "3.460"
Set-TextValue $ws.Cells.Item(17,5) "-0.10%"
Set-TextValue $ws.Cells.Item(17,7) "7"
Set-TextValue $ws.Cells.Item(18,4) "2.381"
Set-TextValue $ws.Cells.Item(18,5) "14.79%"
Set-TextValue $ws.Cells.Item(18,7) "7"
Set-TextValue $ws.Cells.Item(19,4) "0.3428"
Set-TextValue $ws.Cells.Item(19,5) "-1.02%"
Set-TextValue $ws.Cells.Item(19,7) "7"
Set-TextValue $ws.Cells.Item(20,4) "0.1309"
Set-TextValue $ws.Cells.Item(20,5) "-0.01%"
Set-TextValue $ws.Cells.Item(20,7) "7"
Set-TextValue $ws.Cells.Item(21,4) "4.840"
Set-TextValue $ws.Cells.Item(21,5) "1.96%"
Set-TextValue $ws.Cells.Item(21,7) "7"
Set-TextValue $ws.Cells.Item(22,4) "0.2200"
Set-TextValue $ws.Cells.Item(22,5) "-0.19%"
Set-TextValue $ws.Cells.Item(22,7) "7"
Set-TextValue $ws.Cells.Item(23,4) "0.04541"
Set-TextValue $ws.Cells.Item(23,5) "-2.31%"
Set-TextValue $ws.Cells.Item(23,7) "7"
Set-TextValue $ws.Cells.Item(24,4) "0.005090"
Set-TextValue $ws.Cells.Item(24,5) "14.32%"
Set-TextValue $ws.Cells.Item(24,7) "7"
Set-TextValue $ws.Cells.Item(25,4) "0.001234"
Set-TextValue $ws.Cells.Item(25,5) "0.39%"
Set-TextValue $ws.Cells.Item(25,7) "7"
Set-TextValue $ws.Cells.Item(26,5) "7.76%"
Set-TextValue $ws.Cells.Item(26,7) "7"
Set-TextValue $ws.Cells.Item(27,4) "0.0004749"
Set-TextValue $ws.Cells.Item(27,5) "-0.04%"
Set-TextValue $ws.Cells.Item(27,7) "7"
Set-TextValue $ws.Cells.Item(28,7) "7"
Set-TextValue $ws.Cells.Item(29,7) "7"
Set-TextValue $ws.Cells.Item(30,7) "7"
Set-TextValue $ws.Cells.Item(31,7) "7"
Set-TextValue $ws.Cells.Item(32,7) "7"
Set-TextValue $ws.Cells.Item(33,7) "7"
Set-TextValue $ws.Cells.Item(34,7) "7"
Set-TextValue $ws.Cells.Item(35,7) "7"
Set-TextValue $ws.Cells.Item(36,7) "7"
Set-TextValue $ws.Cells.Item(37,7) "7"
Set-TextValue $ws.Cells.Item(38,7) "7"
Set-TextValue $ws.Cells.Item(39,4) "0.01854"
Set-TextValue $ws.Cells.Item(39,5) "4.89%"
Set-TextValue $ws.Cells.Item(39,7) "7"
Set-TextValue $ws.Cells.Item(40,4) "0.04772"
Set-TextValue $ws.Cells.Item(40,5) "0.49%"
Set-TextValue $ws.Cells.Item(40,7) "7"
Set-TextValue $ws.Cells.Item(41,4) "0.007788"
Set-TextValue $ws.Cells.Item(41,5) "-0.46%"
Set-TextValue $ws.Cells.Item(41,7) "7"
Set-TextValue $ws.Cells.Item(42,4) "0.1402"
Set-TextValue $ws.Cells.Item(42,5) "0.71%"
Set-TextValue $ws.Cells.Item(42,7) "7"
Set-TextValue $ws.Cells.Item(43,4) "0.007730"
Set-TextValue $ws.Cells.Item(43,5) "0.93%"
Set-TextValue $ws.Cells.Item(43,7) "7"
Set-TextValue $ws.Cells.Item(44,4) "0.002229"
Set-TextValue $ws.Cells.Item(44,5) "0.05%"
Set-TextValue $ws.Cells.Item(44,7) "7"
Set-TextValue $ws.Cells.Item(45,4) "0.01107"
Set-TextValue $ws.Cells.Item(45,5) "12.54%"
Set-TextValue $ws.Cells.Item(45,7) "7"
Set-TextValue $ws.Cells.Item(46,4) "0.00006414"
Set-TextValue $ws.Cells.Item(46,5) "5.83%"
Set-TextValue $ws.Cells.Item(46,7) "7"
Set-TextValue $ws.Cells.Item(47,5) "-0.03%"
Set-TextValue $ws.Cells.Item(47,7) "7"
Set-TextValue $ws.Cells.Item(48,4) "0.0005800"
Set-TextValue $ws.Cells.Item(48,5) "-0.01%"
Set-TextValue $ws.Cells.Item(48,7) "7"
Set-TextValue $ws.Cells.Item(49,4) "24.50"
Set-TextValue $ws.Cells.Item(49,5) "182.59%"
Set-TextValue $ws.Cells.Item(49,7) "7"
Set-TextValue $ws.Cells.Item(50,4) "0.002000"
Set-TextValue $ws.Cells.Item(50,5) "-25.68%"
Set-TextValue $ws.Cells.Item(50,7) "7"
Set-TextValue $ws.Cells.Item(51,4) "0.00002100"
Set-TextValue $ws.Cells.Item(51,5) "-0.03%"
Set-TextValue $ws.Cells.Item(51,7) "7"
